$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.050840816311689387
$ws.Range("B1").Value = -0.043126609670307188
$ws.Range("A2").Value = -0.059421440580177953
$ws.Range("B2").Value = -0.013588820181848886
$ws.Range("A3").Value = -0.012679583874133631
$ws.Range("B3").Value = -0.035733986014339927
$ws.Range("A4").Value = 0.043744225914486616
$ws.Range("B4").Value = -0.043744225934178115
